$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue 2 4 "63.810.73"
Set-TextValue 2 5 "  +4.80%  "

# Row 3
Set-TextValue 3 4 "2.759.57"
Set-TextValue 3 5 "  +4.41%  "

# Row 4
Set-TextValue 4 4 "1.00"
Set-TextValue 4 5 "  +0.12%  "

# Row 5
Set-TextValue 5 4 "581.48"
Set-TextValue 5 5 "  +0.73%  "

# Row 6
Set-TextValue 6 4 "157.49"
Set-TextValue 6 5 "  +9.41%  "

# Row 7
Set-TextValue 7 2 "XRP"
Set-TextValue 7 3 "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue 7 4 "0.623"
Set-TextValue 7 5 "  +4.28%  "

# Row 8
Set-TextValue 8 2 "USDC"
Set-TextValue 8 3 "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue 8 4 "0.997"
Set-TextValue 8 5 "  -0.09%  "

# Row 9
Set-TextValue 9 4 "2.757.33"
Set-TextValue 9 5 "  +3.75%  "

# Row 10
Set-TextValue 10 4 "6.77"
Set-TextValue 10 5 "  +2.80%  "

# Row 11
Set-TextValue 11 5 "  +3.63%  "

# Row 12
Set-TextValue 12 5 "  +4.14%  "

# Row 13
Set-TextValue 13 5 "  +0.81%  "

# Row 14
Set-TextValue 14 4 "3.240.47"
Set-TextValue 14 5 "  +4.07%  "

# Row 15
Set-TextValue 15 4 "27.53"
Set-TextValue 15 5 "  +4.29%  "

# Row 16
Set-TextValue 16 4 "63.784.28"
Set-TextValue 16 5 "  +4.85%  "

# Row 17
Set-TextValue 17 4 "0.0000155"
Set-TextValue 17 5 "  +7.30%  "

# Row 18
Set-TextValue 18 4 "2.752.20"
Set-TextValue 18 5 "  +3.37%  "

# Row 19
Set-TextValue 19 4 "12.08"
Set-TextValue 19 5 "  +3.84%  "

# Row 20
Set-TextValue 20 4 "4.95"
Set-TextValue 20 5 "  +4.54%  "

# Row 21
Set-TextValue 21 4 "362.87"
Set-TextValue 21 5 "  +3.05%  "

# Row 22
Set-TextValue 22 4 "6.96"
Set-TextValue 22 5 "  +1.09%  "

# Row 23
Set-TextValue 23 4 "0.545"
Set-TextValue 23 5 "  +3.53%  "

# Row 24
Set-TextValue 24 4 "1.00"
Set-TextValue 24 5 "  +0.13%  "

# Row 25
Set-TextValue 25 4 "66.85"
Set-TextValue 25 5 "  +4.51%  "

# Row 26
Set-TextValue 26 5 "  +5.59%  "

# Row 27
Set-TextValue 27 4 "8.65"
Set-TextValue 27 5 "  +3.00%  "

# Row 28
Set-TextValue 28 5 "  -0.15%  "

# Row 29
Set-TextValue 29 4 "0.0₃0927"
Set-TextValue 29 5 "  +14.67%  "

# Row 30
Set-TextValue 30 5 "  +0.95%  "

# Row 31
Set-TextValue 31 4 "7.25"
Set-TextValue 31 5 "  +7.07%  "

# Row 32
Set-TextValue 32 4 "1.29"
Set-TextValue 32 5 "  +19.81%  "

# Row 33
Set-TextValue 33 4 "173.88"
Set-TextValue 33 5 "  +4.35%  "

# Row 34
Set-TextValue 34 5 "  -0.03%  "

# Row 35
Set-TextValue 35 4 "20.64"
Set-TextValue 35 5 "  +3.37%  "

# Row 36
Set-TextValue 36 5 "  +7.67%  "

# Row 37
Set-TextValue 37 5 "  +10.86%  "

# Row 38
Set-TextValue 38 5 "  +8.26%  "

# Row 39
Set-TextValue 39 5 "  +11.40%  "

# Row 40
Set-TextValue 40 4 "4.30"
Set-TextValue 40 5 "  +4.06%  "

# Row 41
Set-TextValue 41 4 "338.28"
Set-TextValue 41 5 "  -1.01%  "

# Row 42
Set-TextValue 42 4 "6.05"
Set-TextValue 42 5 "  +14.70%  "

# Row 43
Set-TextValue 43 4 "39.49"
Set-TextValue 43 5 "  +3.06%  "

# Row 44
Set-TextValue 44 4 "21.93"
Set-TextValue 44 5 "  +7.98%  "

# Row 45
Set-TextValue 45 4 "22.25"
Set-TextValue 45 5 "  +6.84%  "

# Row 46
Set-TextValue 46 5 "  +5.39%  "

# Row 47
Set-TextValue 47 4 "0.649"
Set-TextValue 47 5 "  +4.14%  "

# Row 48
Set-TextValue 48 5 "  +3.72%  "

# Row 49
Set-TextValue 49 4 "138.76"
Set-TextValue 49 5 "  +0.35%  "

# Row 50
Set-TextValue 50 4 "0.103"
Set-TextValue 50 5 "  +3.27%  "

# Row 51
Set-TextValue 51 4 "0.997"
Set-TextValue 51 5 "  -0.09%  "
